$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Relevant Coursework section: originally 3 paragraphs (16, 17, 18) listing
#    courses. Re-grouped into 2 paragraphs with different course ordering,
#    and the leading tab character on the first line is replaced by a
#    first-line indent.
# ---------------------------------------------------------------------------

# Locate the three coursework paragraphs by content (robust to index drift).
$cwStart = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Discrete Structures 1 & 2*") {
        $cwStart = $p.Range.Start
        break
    }
}

$p1 = $d.Paragraphs(1)
$p1 = $null
# Re-fetch paragraph objects from the located start point each time,
# since paragraph indices / objects shift as we edit.
$pA = $d.Range($cwStart, $cwStart).Paragraphs(1)

# Merge the three coursework paragraphs into a single paragraph so we can
# retype its contents cleanly, then split it back into two paragraphs.
$full = $pA.Range
$pmark = $d.Range($full.End - 1, $full.End)
$pmark.Delete()

$pA = $d.Range($cwStart, $cwStart).Paragraphs(1)
$full2 = $pA.Range
$pmark2 = $d.Range($full2.End - 1, $full2.End)
$pmark2.Delete()

# Remove the leading tab character (first character of the merged paragraph).
$pA = $d.Range($cwStart, $cwStart).Paragraphs(1)
$mfull = $pA.Range
$tabRange = $d.Range($mfull.Start, $mfull.Start + 1)
$tabRange.Delete()

# Retype the full merged text with a marker at the point where the
# paragraph should be split back into two.
$pA = $d.Range($cwStart, $cwStart).Paragraphs(1)
$bodyRange = $pA.Range
$newText = "Discrete Structures 1 & 2, Introduction to Java, Database Systems, Operating Systems, Foundations of Computing, |SPLIT|Algorithms and Data Structures, Algorithms II, Continuous Algorithms, Computer Architecture, Software Engineering, Systems Programming, Introduction to Robotics Manipulation, Machine Learning, Introduction to Statistical Natural Language Processing, Deep Question Answering with IBM Watson"
$bodyRange.Text = $newText

# Split the marker into a real paragraph break.
$d.Content.Find.Execute("|SPLIT|", $false, $false, $false, $false, $false, $true, 1, $false, "^p", 2) | Out-Null

$pA = $d.Range($cwStart, $cwStart).Paragraphs(1)
$pA.Format.FirstLineIndent = 36

$pB = $pA.Next()
$pB.Format.LeftIndent = 36

# Move the "_GoBack" bookmark to the end of paragraph B's text (it will be
# relocated there later from its original position near "Working knowledge").
